# Append 12 additional year-data rows (A206:B217) to Sheet1, continuing the
# existing A1:B205 series used as machine-learning X data.
# Column A keeps the same style as the rest of the series (s="1", copied
# from the last existing data row), column B is left with default style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append: index values (col A) and data values (col B)
$newData = @(
  @(204, 0.5247524752475247),
  @(205, 0.6783678367836783),
  @(206, 0.02310231023102316),
  @(207, 0.7711771177117711),
  @(208, 0.4191419141914191),
  @(209, 0.7586044318717585),
  @(210, 0.1551155115511551),
  @(211, 0.2607260726072608),
  @(212, 0.4191419141914191),
  @(213, 0.4191419141914191),
  @(214, 0.4191419141914191),
  @(215, 0.4191419141914191)
)

$startRow = 206
$endRow = $startRow + $newData.Length - 1

# Copy the formatting (style) of the last existing "index" cell (A205) down
# across the new index cells so the appended rows match the existing look.
$ws.Range("A205").Copy()
$ws.Range("A206:A217").PasteSpecial(-4122)  # xlPasteFormats

$r = $startRow
foreach ($pair in $newData) {
  $ws.Cells.Item($r, 1).Value = $pair[0]
  $ws.Cells.Item($r, 2).Value = $pair[1]
  $r = $r + 1
}
